{"js": "// Replace the 25 division-problem answers in the single table of the\n// document. The table has 20 rows x 5 columns; only rows 0, 4, 8, 12, 16\n// (0-based) hold text, one division problem per cell, and every one of\n// those 25 cells gets a new value (row/col positions below match the\n// document's current row/col order exactly).\nconst updates = [\n  [0, 0, \"490\u00f72=245, 0\"],\n  [0, 1, \"941\u00f74=235, 1\"],\n  [0, 2, \"434\u00f78=54, 2\"],\n  [0, 3, \"250\u00f72=125, 0\"],\n  [0, 4, \"297\u00f78=37, 1\"],\n  [4, 0, \"760\u00f78=95, 0\"],\n  [4, 1, \"785\u00f77=112, 1\"],\n  [4, 2, \"918\u00f78=114, 6\"],\n  [4, 3, \"542\u00f78=67, 6\"],\n  [4, 4, \"347\u00f79=38, 5\"],\n  [8, 0, \"755\u00f75=151, 0\"],\n  [8, 1, \"163\u00f76=27, 1\"],\n  [8, 2, \"732\u00f74=183, 0\"],\n  [8, 3, \"632\u00f79=70, 2\"],\n  [8, 4, \"275\u00f77=39, 2\"],\n  [12, 0, \"819\u00f72=409, 1\"],\n  [12, 1, \"755\u00f73=251, 2\"],\n  [12, 2, \"688\u00f76=114, 4\"],\n  [12, 3, \"103\u00f78=12, 7\"],\n  [12, 4, \"475\u00f72=237, 1\"],\n  [16, 0, \"668\u00f74=167, 0\"],\n  [16, 1, \"371\u00f75=74, 1\"],\n  [16, 2, \"893\u00f72=446, 1\"],\n  [16, 3, \"707\u00f74=176, 3\"],\n  [16, 4, \"270\u00f75=54, 0\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, text] of updates) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the single table of the\n# document. The table has 20 rows x 5 columns (1-based for COM); only rows\n# 1, 5, 9, 13, 17 hold text, one division problem per cell, and every one\n# of those 25 cells gets a new value (row/col positions below match the\n# document's current row/col order exactly).\n$updates = @(\n  @(1, 1, \"490\u00f72=245, 0\"),\n  @(1, 2, \"941\u00f74=235, 1\"),\n  @(1, 3, \"434\u00f78=54, 2\"),\n  @(1, 4, \"250\u00f72=125, 0\"),\n  @(1, 5, \"297\u00f78=37, 1\"),\n  @(5, 1, \"760\u00f78=95, 0\"),\n  @(5, 2, \"785\u00f77=112, 1\"),\n  @(5, 3, \"918\u00f78=114, 6\"),\n  @(5, 4, \"542\u00f78=67, 6\"),\n  @(5, 5, \"347\u00f79=38, 5\"),\n  @(9, 1, \"755\u00f75=151, 0\"),\n  @(9, 2, \"163\u00f76=27, 1\"),\n  @(9, 3, \"732\u00f74=183, 0\"),\n  @(9, 4, \"632\u00f79=70, 2\"),\n  @(9, 5, \"275\u00f77=39, 2\"),\n  @(13, 1, \"819\u00f72=409, 1\"),\n  @(13, 2, \"755\u00f73=251, 2\"),\n  @(13, 3, \"688\u00f76=114, 4\"),\n  @(13, 4, \"103\u00f78=12, 7\"),\n  @(13, 5, \"475\u00f72=237, 1\"),\n  @(17, 1, \"668\u00f74=167, 0\"),\n  @(17, 2, \"371\u00f75=74, 1\"),\n  @(17, 3, \"893\u00f72=446, 1\"),\n  @(17, 4, \"707\u00f74=176, 3\"),\n  @(17, 5, \"270\u00f75=54, 0\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($u in $updates) {\n  $row = $u[0]\n  $col = $u[1]\n  $text = $u[2]\n  $t.Cell($row, $col).Range.Text = $text\n}\n"}
